$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates as described by the commit diff (updated crypto price/volume data,
# including two pairs of rows whose entire record order was swapped: rows 29/30 and 46/47).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.569.13'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.93%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.996.93'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.83%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '382.80'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.49'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +2.19%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +1.95%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.594'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.92'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.65%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0861'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.467.59'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.82'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +3.40%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '18.45'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.20%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.987.92'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '11.13'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +3.82%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.999'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '51.583.02'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.09'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0962'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.51'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.73%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '267.96'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.21'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.87%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.88'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.77%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.46'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.09%  '
$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.167'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.91%  '
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '26.08'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.109'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.35'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +3.33%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.75'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +4.17%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '51.59'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.79%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.04'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '16.78'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +3.49%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +3.29%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.96%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '125.00'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +4.05%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.65'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +9.32%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.55'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.03'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.38'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +3.44%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.271'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.046.44'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0335'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +3.38%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +17.63%  '
